$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use an existing date-formatted cell as the style template so the same
# cellXf (s="1") gets reused instead of a brand-new one being created.
$dateTemplate = $ws.Range("A108")

# Row 109
$ws.Range("A109").Value = 45498.2916666667
$dateTemplate.Copy()
$ws.Range("A109").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("B109").Value = 0
$ws.Range("C109").Value = 3.14000010490417
$ws.Range("D109").Value = 3.14000010490417
$ws.Range("E109").Value = 3.14000010490417
$ws.Range("F109").Value = 3.14000010490417
$ws.Range("G109").Formula = "'3.14000010490417"
$ws.Range("G109").Style = "Normal"
$ws.Range("H109").Value = "ESPE.MI"

# Row 110
$ws.Range("A110").Value = 45499.3369907407
$dateTemplate.Copy()
$ws.Range("A110").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("B110").Value = 5250
$ws.Range("C110").Value = 3.14000010490417
$ws.Range("D110").Value = 3.04999995231628
$ws.Range("E110").Value = 3.04999995231628
$ws.Range("F110").Value = 3.14000010490417
$ws.Range("G110").Formula = "'3.14000010490417"
$ws.Range("G110").Style = "Normal"
$ws.Range("H110").Value = "ESPE.MI"
